$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.780.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.222.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.14"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.103"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.85"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.563.72"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.830"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.236.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.763.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000104"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.04"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.85%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.24"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.08"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.39"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0791"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.51"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.66"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.24"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0301"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.97"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.10"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.61"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.195"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.65"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.22"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0999"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.68"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.441.44"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.08%  "
